$d = $word.ActiveDocument

# Locate "QUINTA" (CLAUSULA QUINTA heading) and append a new run
# containing ":" right after it, matching the bold Arial Rounded MT
# Bold / size-20-halfpoints formatting used by the clause heading.
$rng = $d.Content
$rng.Find.Execute("QUINTA", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null

$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertAfter(":")
$insertPoint.Font.Name = "Arial Rounded MT Bold"
$insertPoint.Font.Bold = $true
$insertPoint.Font.BoldBi = $true
$insertPoint.Font.Size = 10
$insertPoint.Font.SizeBi = 10
